$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 4 ("Douchebag" entry) with the new "Afterthought" entry
$ws.Range("A4").Value = "Afterthought"
$ws.Range("B4").Value = "পরে আসা চিন্তা"
$ws.Range("C4").Value = "As an afterthought, he gave me his phone number"

# Remove the old row 5 ("Meaning" / "অর্থ ") entirely
$ws.Rows.Item(5).Delete()

# Update the active selection to match the saved view state
$ws.Range("B5").Select()
